$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCredentials")
$ws.Range("C3:D3").Borders.LineStyle = -4142
$ws.Range("C3:D3").Interior.Pattern = -4142
$ws.Range("C3").Value = "authorprod@knowledgehut.com"
$ws.Range("D3").Value = "Password@123"
